{"js": "const body = context.document.body;\n\nconst replacements = [\n  [\"569\u00f74=142, 1\", \"508\u00f78=63, 4\"],\n  [\"965\u00f77=137, 6\", \"469\u00f75=93, 4\"],\n  [\"449\u00f77=64, 1\", \"446\u00f78=55, 6\"],\n  [\"434\u00f74=108, 2\", \"259\u00f77=37, 0\"],\n  [\"268\u00f78=33, 4\", \"505\u00f74=126, 1\"],\n  [\"285\u00f74=71, 1\", \"652\u00f77=93, 1\"],\n  [\"999\u00f76=166, 3\", \"877\u00f73=292, 1\"],\n  [\"435\u00f79=48, 3\", \"926\u00f76=154, 2\"],\n  [\"981\u00f79=109, 0\", \"785\u00f77=112, 1\"],\n  [\"907\u00f72=453, 1\", \"644\u00f76=107, 2\"],\n  [\"520\u00f73=173, 1\", \"476\u00f73=158, 2\"],\n  [\"105\u00f76=17, 3\", \"766\u00f75=153, 1\"],\n  [\"990\u00f78=123, 6\", \"288\u00f74=72, 0\"],\n  [\"172\u00f79=19, 1\", \"843\u00f74=210, 3\"],\n  [\"955\u00f79=106, 1\", \"498\u00f72=249, 0\"],\n  [\"791\u00f75=158, 1\", \"220\u00f74=55, 0\"],\n  [\"551\u00f76=91, 5\", \"323\u00f72=161, 1\"],\n  [\"853\u00f73=284, 1\", \"821\u00f78=102, 5\"],\n  [\"702\u00f75=140, 2\", \"182\u00f74=45, 2\"],\n  [\"552\u00f73=184, 0\", \"698\u00f74=174, 2\"],\n  [\"259\u00f73=86, 1\", \"505\u00f77=72, 1\"],\n  [\"697\u00f72=348, 1\", \"948\u00f74=237, 0\"],\n  [\"195\u00f77=27, 6\", \"925\u00f79=102, 7\"],\n  [\"861\u00f79=95, 6\", \"632\u00f74=158, 0\"],\n  [\"128\u00f79=14, 2\", \"918\u00f73=306, 0\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"569\u00f74=142, 1\", \"508\u00f78=63, 4\"),\n    @(\"965\u00f77=137, 6\", \"469\u00f75=93, 4\"),\n    @(\"449\u00f77=64, 1\", \"446\u00f78=55, 6\"),\n    @(\"434\u00f74=108, 2\", \"259\u00f77=37, 0\"),\n    @(\"268\u00f78=33, 4\", \"505\u00f74=126, 1\"),\n    @(\"285\u00f74=71, 1\", \"652\u00f77=93, 1\"),\n    @(\"999\u00f76=166, 3\", \"877\u00f73=292, 1\"),\n    @(\"435\u00f79=48, 3\", \"926\u00f76=154, 2\"),\n    @(\"981\u00f79=109, 0\", \"785\u00f77=112, 1\"),\n    @(\"907\u00f72=453, 1\", \"644\u00f76=107, 2\"),\n    @(\"520\u00f73=173, 1\", \"476\u00f73=158, 2\"),\n    @(\"105\u00f76=17, 3\", \"766\u00f75=153, 1\"),\n    @(\"990\u00f78=123, 6\", \"288\u00f74=72, 0\"),\n    @(\"172\u00f79=19, 1\", \"843\u00f74=210, 3\"),\n    @(\"955\u00f79=106, 1\", \"498\u00f72=249, 0\"),\n    @(\"791\u00f75=158, 1\", \"220\u00f74=55, 0\"),\n    @(\"551\u00f76=91, 5\", \"323\u00f72=161, 1\"),\n    @(\"853\u00f73=284, 1\", \"821\u00f78=102, 5\"),\n    @(\"702\u00f75=140, 2\", \"182\u00f74=45, 2\"),\n    @(\"552\u00f73=184, 0\", \"698\u00f74=174, 2\"),\n    @(\"259\u00f73=86, 1\", \"505\u00f77=72, 1\"),\n    @(\"697\u00f72=348, 1\", \"948\u00f74=237, 0\"),\n    @(\"195\u00f77=27, 6\", \"925\u00f79=102, 7\"),\n    @(\"861\u00f79=95, 6\", \"632\u00f74=158, 0\"),\n    @(\"128\u00f79=14, 2\", \"918\u00f73=306, 0\"),\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2) | Out-Null\n}\n"}
